$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table 4")
$ws.Range("A1").Value = "TEST"
